$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.085.20"
$ws.Range("D3").Value = "2.639.83"
$ws.Range("E3").Value = "  +10.73%  "
$ws.Range("D5").Value = "313.22"
$ws.Range("E5").Value = "  +6.87%  "
$ws.Range("D6").Value = "104.46"
$ws.Range("E6").Value = "  +10.77%  "
$ws.Range("D7").Value = "0.612"
$ws.Range("E7").Value = "  +10.39%  "
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "0.596"
$ws.Range("E9").Value = "  +19.94%  "
$ws.Range("D10").Value = "39.66"
$ws.Range("E10").Value = "  +16.96%  "
$ws.Range("D11").Value = "55.32"
$ws.Range("E11").Value = "  +3.47%  "
$ws.Range("D12").Value = "0.0855"
$ws.Range("E12").Value = "  +10.62%  "
$ws.Range("D13").Value = "8.43"
$ws.Range("E13").Value = "  +21.86%  "
$ws.Range("D14").Value = "3.041.77"
$ws.Range("E14").Value = "  +10.77%  "
$ws.Range("E15").Value = "  +3.30%  "
$ws.Range("D16").Value = "2.661.06"
$ws.Range("E16").Value = "  +11.94%  "
$ws.Range("D17").Value = "0.944"
$ws.Range("E17").Value = "  +15.31%  "
$ws.Range("D18").Value = "15.37"
$ws.Range("E18").Value = "  +10.61%  "
$ws.Range("D19").Value = "47.577.28"
$ws.Range("E19").Value = "  +5.59%  "
$ws.Range("E20").Value = "  +11.79%  "
$ws.Range("D21").Value = "13.25"
$ws.Range("E21").Value = "  +7.54%  "
$ws.Range("E22").Value = "  +12.16%  "
$ws.Range("D23").Value = "72.57"
$ws.Range("E23").Value = "  +9.84%  "
$ws.Range("D24").Value = "272.03"
$ws.Range("E24").Value = "  +14.34%  "
$ws.Range("E25").Value = "  +13.29%  "
$ws.Range("D26").Value = "2.22"
$ws.Range("E26").Value = "  +18.70%  "
$ws.Range("D27").Value = "30.37"
$ws.Range("E27").Value = "  +45.23%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "10.74"
$ws.Range("E30").Value = "  +13.23%  "
$ws.Range("D31").Value = "40.20"
$ws.Range("E31").Value = "  +8.14%  "
$ws.Range("E32").Value = "  +4.97%  "
$ws.Range("D33").Value = "6.22"
$ws.Range("E33").Value = "  +15.66%  "
$ws.Range("D34").Value = "3.73"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "2.27"
$ws.Range("E35").Value = "  +16.77%  "
$ws.Range("D36").Value = "0.0857"
$ws.Range("E36").Value = "  +13.45%  "
$ws.Range("E37").Value = "  +6.20%  "
$ws.Range("D38").Value = "151.92"
$ws.Range("E38").Value = "  +3.44%  "
$ws.Range("E39").Value = "  +11.35%  "
$ws.Range("E40").Value = "  +10.10%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "16.66"
$ws.Range("E41").Value = "  +13.56%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "23.24"
$ws.Range("E42").Value = "  +55.54%  "
$ws.Range("D43").Value = "4.29"
$ws.Range("E43").Value = "  +15.71%  "
$ws.Range("D44").Value = "3.74"
$ws.Range("E44").Value = "  +18.61%  "
$ws.Range("E45").Value = "  +14.23%  "
$ws.Range("D46").Value = "2.188.69"
$ws.Range("E46").Value = "  +12.01%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D47").Value = "95.82"
$ws.Range("E47").Value = "  +7.71%  "
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("D49").Value = "10.07"
$ws.Range("E49").Value = "  +19.82%  "
$ws.Range("D50").Value = "114.73"
$ws.Range("E50").Value = "  +15.77%  "
$ws.Range("E51").Value = "  +7.30%  "
